# Adding checks that every container has a process.
# -> Add an "outlet"/"time" log (discharge column L, date-serial column N)
#    to the glacier snowpack model sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells on row 6.
# NB: write "time" before "outlet" so the shared-strings table gets the
# same ordering as the target workbook (time=13, outlet=14).
$ws.Range("N6").Value = "time"
$ws.Range("L6").Value = "outlet"

# New "time" column (N7:N16) - sequential day-serial numbers.
$ws.Range("N7").Value  = 58849
$ws.Range("N8").Value  = 58850
$ws.Range("N9").Value  = 58851
$ws.Range("N10").Value = 58852
$ws.Range("N11").Value = 58853
$ws.Range("N12").Value = 58854
$ws.Range("N13").Value = 58855
$ws.Range("N14").Value = 58856
$ws.Range("N15").Value = 58857
$ws.Range("N16").Value = 58858

# L12 was re-typed directly (breaking it out of the L8:L16 shared formula)
# even though the result is identical.
$ws.Range("L12").Formula = '=F12+H12+J12*$M$2'

# Leave the selection on the new header cell, matching the saved view state.
$ws.Range("L6").Select()
